$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.050.69"
$ws.Range("E2").Value = "  +5.06%  "
$ws.Range("D3").Value = "3.520.27"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.00"
$ws.Range("E5").Value = "  +4.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.01"
$ws.Range("E6").Value = "  +6.99%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.521.42"
$ws.Range("E8").Value = "  +2.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.125"
$ws.Range("E11").Value = "  +5.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("E12").Value = "  +4.31%  "
$ws.Range("D13").Value = "4.125.15"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.22"
$ws.Range("E15").Value = "  +4.22%  "
$ws.Range("E16").Value = "  +4.44%  "
$ws.Range("D17").Value = "67.022.41"
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("D18").Value = "3.522.47"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("E19").Value = "  +4.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.05"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.42"
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.98"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.60"
$ws.Range("E23").Value = "  +3.26%  "
$ws.Range("E24").Value = "  +10.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.21"
$ws.Range("E27").Value = "  +5.49%  "
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.43"
$ws.Range("E30").Value = "  +5.69%  "
$ws.Range("E31").Value = "  +6.23%  "
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.63"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.47"
$ws.Range("E34").Value = "  +7.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.59"
$ws.Range("E36").Value = "  +5.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.22"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.902"
$ws.Range("E38").Value = "  +6.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").Value = "  +5.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0754"
$ws.Range("E40").Value = "  +4.29%  "
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.59"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.72"
$ws.Range("E43").Value = "  +5.05%  "
$ws.Range("D44").Value = "2.840.76"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.51"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.42"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("E47").Value = "  +7.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0314"
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "353.69"
$ws.Range("E49").Value = "  +6.12%  "
$ws.Range("E50").Value = "  +4.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.62"
$ws.Range("E51").Value = "  +12.09%  "
